$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33, shifting existing rows 33..71 down to 34..72.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly data point.
$ws.Cells.Item(33, 1).Value = 4
$ws.Cells.Item(33, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(33, 3).Value = "Los Lagos"
$ws.Cells.Item(33, 4).Value = 44495
$ws.Cells.Item(33, 5).Value = 10
$ws.Cells.Item(33, 6).Value = 100112052
$ws.Cells.Item(33, 7).Value = "Albahaca"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 90
$ws.Cells.Item(33, 11).Value = 5000
$ws.Cells.Item(33, 12).Value = 5000
$ws.Cells.Item(33, 13).Value = 5000
$ws.Cells.Item(33, 14).Value = "$/paquete"
$ws.Cells.Item(33, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value = 5000
$ws.Cells.Item(33, 17).Value = 1
$ws.Cells.Item(33, 18).Value = "Hortaliza"
